$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sub Contract works")

# --- Section headers become bold ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A6").Font.Bold = $true
$ws.Range("A9").Font.Bold = $true

# --- "Erathwork supports" / m2 line (row 15) is removed entirely ---
$ws.Range("A15:B15").EntireRow.Delete()

# --- Update unit-of-measure values for the "Excavation & Erathwork" items ---
$ws.Range("B10").Value = "item"
$ws.Range("B11").Value = "item"
$ws.Range("B12").Value = "nr"
$ws.Range("B13").Value = "cube"
$ws.Range("B14").Value = "cube"
$ws.Range("B15").Value = "item"
$ws.Range("B16").Value = "cube"

# --- Collapse the 3 blank rows that used to separate this block from "Concrete works" ---
$ws.Range("A18:B20").EntireRow.Delete()

# --- Insert the new "Concrete works" line items + trailing blank separator row ---
# (done before bolding the header row so the new rows don't inherit its bold font)
$ws.Range("A19:B23").EntireRow.Insert()

$ws.Range("A19").Value = "Screed concrete"
$ws.Range("B19").Value = "ft2"
$ws.Range("A20").Value = "Form work"
$ws.Range("B20").Value = "ft2"
$ws.Range("A21").Value = "Reinforcement"
$ws.Range("B21").Value = "kg"
$ws.Range("A22").Value = "slab concrete"
$ws.Range("B22").Value = "bag"

# --- "Concrete works" header (row 18) becomes bold ---
$ws.Range("A18").Font.Bold = $true

# --- Insert 6 new blank rows under the "Masonary works" header (row 24) ---
# (done before bolding the header row so the new rows stay independently styled)
$ws.Range("A25:B30").EntireRow.Insert()

# --- "Masonary works" header (row 24) becomes bold ---
$ws.Range("A24").Font.Bold = $true

# --- The 6 new blank rows underneath it are bold too ---
$ws.Range("A25:A30").Font.Bold = $true

# --- Restore the original page setup / view state ---
$ws.PageSetup.Orientation = 1
$ws.Rows("25:25").Select()
